$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Mark additional rows as "done" (Good style), matching already-styled rows.
$ws1.Range("B9").Style = "Good"
$ws1.Range("B20").Style = "Good"
$ws1.Range("B21").Style = "Good"

# Update selection on the original sheet.
$ws1.Range("B7:B21").Select()

# Add a new worksheet at the end of the workbook with interactor instructions.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Range("A1").Value = "ctrl + x"
$ws2.Range("B1").Value = "interactor instructions"
$ws2.Range("A2").Value = "shift + c "
$ws2.Range("B2").Value = "polygong offset"
$ws2.Range("A3").Value = "use buffers for nodes and edges"
$ws2.Range("A4").Select()
